$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp title in A1
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 20:46"

# Full country table, re-sorted descending by Casos totales, with newly reported
# countries (Siria, Granada) inserted and updated case counts for all affected rows.
$data = @(
  @(4, "China", 81054, 46, 72440, 5353, 1845, 6, 3261),
  @(5, "Italia", 59138, 5560, 7024, 46638, 3000, 651, 5476),
  @(6, "Estados Unidos", 32356, 8149, 178, 31764, 795, 112, 414),
  @(7, "España", 28603, 3107, 2125, 24722, 1785, 375, 1756),
  @(8, "Alemania", 24806, 2442, 266, 24447, 2, 9, 93),
  @(9, "Iran", 21638, 1028, 7635, 12318, 0, 129, 1685),
  @(10, "Francia", 16018, 1559, 2200, 13144, 1746, 112, 674),
  @(11, "Corea del Sur", 8897, 98, 2909, 5884, 59, 2, 104),
  @(12, "Suiza", 7474, 611, 131, 7245, 141, 18, 98),
  @(13, "Reino Unido", 5683, 665, 93, 5309, 20, 48, 281),
  @(14, "Paises Bajos", 4204, 573, 2, 4023, 354, 43, 179),
  @(15, "Belgica", 3401, 586, 263, 3063, 288, 8, 75),
  @(16, "Austria", 3302, 310, 9, 3277, 15, 8, 16),
  @(17, "Noruega", 2263, 99, 6, 2250, 28, 0, 7),
  @(18, "Suecia", 1931, 161, 16, 1894, 68, 1, 21),
  @(19, "Portugal", 1600, 320, 5, 1581, 26, 2, 14),
  @(20, "Canada", 1426, 98, 14, 1392, 1, 1, 20),
  @(21, "Dinamarca", 1395, 69, 1, 1381, 42, 0, 13),
  @(22, "Australia", 1353, 281, 46, 1300, 2, 0, 7),
  @(23, "Malasia", 1306, 123, 139, 1157, 26, 2, 10),
  @(24, "Turquia", 1236, 289, 0, 1206, 0, 9, 30),
  @(25, "Brasil", 1209, 31, 2, 1189, 18, 0, 18),
  @(26, "Chequia", 1120, 125, 6, 1113, 19, 1, 1),
  @(27, "Japon", 1086, 32, 235, 815, 57, 0, 36),
  @(28, "Israel", 1071, 188, 37, 1033, 18, 0, 1),
  @(29, "Irlanda", 906, 121, 5, 897, 29, 1, 4),
  @(30, "Luxemburgo", 798, 128, 6, 784, 3, 0, 8),
  @(31, "Ecuador", 789, 257, 3, 772, 2, 7, 14),
  @(32, "Crucero", 712, 0, 567, 137, 15, 0, 8),
  @(33, "Pakistan", 646, 1, 13, 629, 0, 1, 4),
  @(34, "Chile", 632, 95, 8, 623, 7, 0, 1),
  @(35, "Polonia", 627, 91, 13, 607, 3, 2, 7),
  @(36, "Finlandia", 626, 103, 10, 615, 12, 0, 1),
  @(37, "Grecia", 624, 94, 19, 590, 18, 2, 15),
  @(38, "Tailandia", 599, 188, 44, 554, 7, 0, 1),
  @(39, "Islandia", 568, 95, 5, 562, 1, 0, 1),
  @(40, "Indonesia", 514, 64, 29, 437, 0, 10, 48),
  @(41, "Arabia Saudita", 511, 119, 17, 494, 0, 0, 0),
  @(42, "Catar", 481, 0, 27, 454, 6, 0, 0),
  @(43, "Singapur", 455, 23, 144, 309, 14, 0, 2),
  @(44, "Rumania", 433, 66, 64, 367, 14, 2, 2),
  @(45, "Eslovenia", 414, 31, 0, 412, 12, 1, 2),
  @(46, "India", 396, 64, 24, 365, 0, 2, 7),
  @(47, "Filipinas", 380, 73, 15, 340, 1, 6, 25),
  @(48, "Rusia", 367, 61, 16, 350, 0, 0, 1),
  @(49, "Peru", 363, 45, 1, 357, 5, 0, 5),
  @(50, "Barein", 332, 22, 149, 181, 4, 1, 2),
  @(51, "Estonia", 326, 20, 4, 322, 2, 0, 0),
  @(52, "Hong Kong", 317, 43, 100, 213, 4, 0, 4),
  @(53, "Egipto", 294, 0, 42, 242, 0, 0, 10),
  @(54, "Sudafrica", 274, 34, 2, 272, 0, 0, 0),
  @(55, "Croacia", 254, 48, 5, 248, 5, 0, 1),
  @(56, "Mexico", 251, 48, 4, 245, 1, 0, 2),
  @(57, "Libano", 248, 18, 8, 236, 4, 0, 4),
  @(58, "Panama", 245, 0, 1, 241, 7, 0, 3),
  @(59, "Irak", 233, 19, 57, 156, 0, 3, 20),
  @(60, "Colombia", 231, 35, 3, 226, 0, 2, 2),
  @(61, "Argentina", 225, 67, 27, 194, 0, 0, 4),
  @(62, "Serbia", 222, 51, 2, 218, 4, 1, 2),
  @(63, "Republica Dominicana", 202, 90, 0, 199, 0, 0, 3),
  @(64, "Argelia", 201, 62, 65, 119, 0, 2, 17),
  @(65, "Armenia", 194, 34, 2, 192, 6, 0, 0),
  @(66, "Kuwait", 188, 12, 30, 158, 5, 0, 0),
  @(67, "Bulgaria", 185, 22, 3, 179, 3, 0, 3),
  @(68, "Eslovaquia", 185, 7, 7, 178, 2, 0, 0),
  @(69, "San Marino", 175, 15, 4, 151, 13, 0, 20),
  @(70, "Taiwan", 169, 16, 28, 139, 0, 0, 2),
  @(71, "Emiratos Arabes Unidos", 153, 0, 38, 113, 2, 0, 2),
  @(72, "Letonia", 139, 15, 1, 138, 0, 0, 0),
  @(73, "Uruguay", 135, 25, 0, 135, 2, 0, 0),
  @(74, "Hungria", 131, 28, 16, 109, 6, 2, 6),
  @(75, "Lituania", 129, 30, 1, 127, 1, 0, 1),
  @(76, "Costa Rica", 117, 0, 2, 113, 2, 0, 2),
  @(77, "Republica de Macedonia", 115, 30, 1, 113, 1, 1, 1),
  @(78, "Islas Feroe", 115, 23, 3, 112, 0, 0, 0),
  @(79, "Principado de Andorra", 113, 25, 1, 111, 2, 1, 1),
  @(80, "Vietnam", 113, 19, 17, 96, 2, 0, 0),
  @(81, "Jordania", 112, 12, 1, 111, 0, 0, 0),
  @(82, "Marruecos", 109, 13, 3, 103, 1, 0, 3),
  @(83, "Republica de Chipre", 95, 11, 3, 91, 3, 0, 1),
  @(84, "Moldavia", 94, 14, 1, 92, 3, 0, 1),
  @(85, "Bosnia y Herzegovina", 94, 1, 2, 91, 1, 0, 1),
  @(86, "Malta", 90, 17, 2, 88, 1, 0, 0),
  @(87, "Albania", 89, 13, 2, 85, 2, 0, 2),
  @(88, "Brunei", 88, 5, 2, 86, 2, 0, 0),
  @(89, "Camboya", 84, 31, 2, 82, 0, 0, 0),
  @(90, "Sri Lanka", 82, 5, 3, 79, 2, 0, 0),
  @(91, "Bielorrusia", 76, 0, 15, 61, 0, 0, 0),
  @(92, "Tunez", 75, 15, 1, 71, 7, 2, 3),
  @(93, "Burkina Faso", 75, 11, 5, 66, 0, 1, 4),
  @(94, "Venezuela", 70, 0, 15, 55, 2, 0, 0),
  @(95, "Senegal", 67, 11, 5, 62, 0, 0, 0),
  @(96, "Nueva Zelanda", 66, 14, 0, 66, 0, 0, 0),
  @(97, "Azerbaiyan", 65, 12, 11, 53, 0, 0, 1),
  @(98, "Kazajistan", 59, 5, 0, 59, 0, 0, 0),
  @(99, "Estado de Palestina", 59, 6, 17, 42, 0, 0, 0),
  @(100, "Guadalupe", 56, 0, 0, 55, 4, 0, 1),
  @(101, "Oman", 55, 3, 17, 38, 0, 0, 0),
  @(102, "Georgia", 54, 5, 1, 53, 1, 0, 0),
  @(103, "Trinidad yTobago", 50, 1, 0, 50, 0, 0, 0),
  @(104, "Reunion", 47, 0, 1, 46, 0, 0, 0),
  @(105, "Ucrania", 47, 0, 1, 43, 0, 0, 3),
  @(106, "Uzbekistan", 43, 2, 0, 43, 0, 0, 0),
  @(107, "Camerun", 40, 0, 2, 38, 0, 0, 0),
  @(108, "Afganistan", 40, 16, 1, 38, 0, 1, 1),
  @(109, "Liechtenstein", 37, 0, 0, 37, 0, 0, 0),
  @(110, "Martinica", 37, 0, 0, 36, 7, 0, 1),
  @(111, "Cuba", 35, 14, 0, 34, 0, 0, 1),
  @(112, "Consejo Danes para los Refugiados", 30, 7, 0, 29, 0, 0, 1),
  @(113, "Nigeria", 30, 8, 2, 28, 0, 0, 0),
  @(114, "Guam", 27, 12, 0, 26, 0, 1, 1),
  @(115, "Banglades", 27, 3, 3, 22, 0, 0, 2),
  @(116, "Honduras", 26, 2, 0, 26, 0, 0, 0),
  @(117, "Bolivia", 24, 5, 0, 24, 0, 0, 0),
  @(118, "Mauricio", 24, 10, 0, 22, 1, 1, 2),
  @(119, "Monaco", 23, 5, 1, 22, 0, 0, 0),
  @(120, "Puerto Rico", 23, 2, 0, 22, 0, 0, 1),
  @(121, "Paraguay", 22, 0, 0, 21, 1, 0, 1),
  @(122, "Macao", 22, 3, 10, 12, 0, 0, 0),
  @(123, "Montenegro", 21, 5, 0, 21, 0, 0, 0),
  @(124, "Ghana", 21, 0, 0, 20, 0, 0, 1),
  @(125, "Jamaica", 19, 0, 2, 16, 0, 0, 1),
  @(126, "Guayana Francesa", 18, 0, 0, 18, 0, 0, 0),
  @(127, "Guyana", 18, 0, 0, 17, 0, 0, 1),
  @(128, "Ruanda", 17, 0, 0, 17, 0, 0, 0),
  @(129, "Guatemala", 17, 0, 0, 16, 0, 0, 1),
  @(130, "Togo", 16, 0, 0, 16, 0, 0, 0),
  @(131, "Polinesia Francesa", 15, 0, 0, 15, 0, 0, 0),
  @(132, "Kenia", 15, 8, 0, 15, 0, 0, 0),
  @(133, "Gibraltar", 15, 5, 2, 13, 0, 0, 0),
  @(134, "Kirguistan", 14, 0, 0, 14, 0, 0, 0),
  @(135, "Barbados", 14, 0, 0, 14, 0, 0, 0),
  @(136, "Costa de Marfil", 14, 0, 1, 13, 0, 0, 0),
  @(137, "Maldivas", 13, 0, 3, 10, 0, 0, 0),
  @(138, "Tanzania", 12, 6, 0, 12, 0, 0, 0),
  @(139, "Mayotte", 11, 0, 0, 11, 0, 0, 0),
  @(140, "Etiopia", 11, 2, 0, 11, 0, 0, 0),
  @(141, "Mongolia", 10, 0, 0, 10, 0, 0, 0),
  @(142, "Aruba", 8, 3, 1, 7, 0, 0, 0),
  @(143, "Seychelles", 7, 0, 0, 7, 0, 0, 0),
  @(144, "Islas Virgenes de los Estados Unidos", 6, 0, 0, 6, 0, 0, 0),
  @(145, "Guinea Ecuatorial", 6, 0, 0, 6, 0, 0, 0),
  @(146, "Isla de Man", 5, 3, 0, 5, 0, 0, 0),
  @(147, "San Martin (Parte Francesa)", 5, 0, 0, 5, 0, 0, 0),
  @(148, "Surinam", 5, 0, 0, 5, 0, 0, 0),
  @(149, "Gabon", 5, 0, 0, 4, 0, 0, 1),
  @(150, "Suazilandia", 4, 3, 0, 4, 0, 0, 0),
  @(151, "Nueva Caledonia", 4, 0, 0, 4, 0, 0, 0),
  @(152, "Bahamas", 4, 0, 0, 4, 0, 0, 0),
  @(153, "Zambia", 3, 1, 0, 3, 0, 0, 0),
  @(154, "El Salvador", 3, 0, 0, 3, 0, 0, 0),
  @(155, "Liberia", 3, 0, 0, 3, 0, 0, 0),
  @(156, "Congo", 3, 0, 0, 3, 0, 0, 0),
  @(157, "Namibia", 3, 0, 0, 3, 0, 0, 0),
  @(158, "Madagascar", 3, 0, 0, 3, 0, 0, 0),
  @(159, "Republica de Africa Central", 3, 0, 0, 3, 0, 0, 0),
  @(160, "Zimbabue", 3, 0, 0, 3, 0, 0, 0),
  @(161, "San Bartolome", 3, 0, 0, 3, 0, 0, 0),
  @(162, "Cabo Verde", 3, 0, 0, 3, 0, 0, 0),
  @(163, "Islas Caimanes", 3, 0, 0, 2, 0, 0, 1),
  @(164, "Curazao", 3, 0, 0, 2, 0, 0, 1),
  @(165, "Haiti", 2, 0, 0, 2, 0, 0, 0),
  @(166, "Groenlandia", 2, 0, 0, 2, 0, 0, 0),
  @(167, "Niger", 2, 1, 0, 2, 0, 0, 0),
  @(168, "Angola", 2, 0, 0, 2, 0, 0, 0),
  @(169, "Butan", 2, 0, 0, 2, 0, 0, 0),
  @(170, "Bermudas", 2, 0, 0, 2, 0, 0, 0),
  @(171, "Mauritania", 2, 0, 0, 2, 0, 0, 0),
  @(172, "Fiyi", 2, 0, 0, 2, 0, 0, 0),
  @(173, "Santa Lucia", 2, 0, 0, 2, 0, 0, 0),
  @(174, "Nicaragua", 2, 0, 0, 2, 0, 0, 0),
  @(175, "Guinea", 2, 0, 0, 2, 0, 0, 0),
  @(176, "Benin", 2, 0, 0, 2, 0, 0, 0),
  @(177, "Sudan", 2, 0, 0, 1, 0, 0, 1),
  @(178, "Siria", 1, 1, 0, 1, 0, 0, 0),
  @(179, "Granada", 1, 1, 0, 1, 0, 0, 0),
  @(180, "Papua Nueva Guinea", 1, 0, 0, 1, 0, 0, 0),
  @(181, "Timor Oriental", 1, 0, 0, 1, 0, 0, 0),
  @(182, "Eritrea", 1, 0, 0, 1, 0, 0, 0),
  @(183, "San Martin (Parte Holandesa)", 1, 0, 0, 1, 0, 0, 0),
  @(184, "Montserrat", 1, 0, 0, 1, 0, 0, 0),
  @(185, "Uganda", 1, 0, 0, 1, 0, 0, 0),
  @(186, "Gambia", 1, 0, 0, 1, 0, 0, 0),
  @(187, "Republica de Yibuti", 1, 0, 0, 1, 0, 0, 0),
  @(188, "Santa Sede", 1, 0, 0, 1, 0, 0, 0),
  @(189, "San Vicente y las Granadinas", 1, 0, 0, 1, 0, 0, 0),
  @(190, "Mozambique", 1, 1, 0, 1, 0, 0, 0),
  @(191, "Republica del Chad", 1, 0, 0, 1, 0, 0, 0),
  @(192, "Somalia", 1, 0, 0, 1, 0, 0, 0),
  @(193, "Antigua y Barbuda", 1, 0, 0, 1, 0, 0, 0),
  @(194, "Nepal", 1, 0, 1, 0, 0, 0, 0)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
    $ws.Cells.Item($r, 5).Value = $row[5]
    $ws.Cells.Item($r, 6).Value = $row[6]
    $ws.Cells.Item($r, 7).Value = $row[7]
    $ws.Cells.Item($r, 8).Value = $row[8]
}

"done: wrote " + $data.Count + " rows"